$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow columns A:C from 39 chars wide down to ~36.57 chars wide ---
$ws.Range("A1:C1").ColumnWidth = 35.67

# --- Reset the view: scroll back so column A is visible again and collapse
#     the lingering selection back down to A1 (the saved file had scrolled
#     to topLeftCell="C1" with U4 selected) ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A1").Select() | Out-Null

# --- Add the new "2023" column (T), mirroring the formatting already used
#     for the "2022" column (S) for both the year header row and the data
#     row ---
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("T5").Value = 40

$excel.CutCopyMode = $false
